$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.239.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "'1.849.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'245.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").Value = "'0.7003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "'0.9996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.07730"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'0.3068"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").Value = "'23.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("D11").Value = "'0.07823"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").Value = "'92.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.849.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.138"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D15").Value = "'0.6868"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").Value = "'6.629"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.96%  "

$ws.Range("D17").Value = "'0.000008331"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").Value = "'29.191.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").Value = "'242.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.08%  "

$ws.Range("D20").Value = "'2.089.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.22%  "

$ws.Range("D21").Value = "'12.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").Value = "'7.527"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'0.9996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "'0.1513"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("D26").Value = "'159.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").Value = "'8.836"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("D29").Value = "'1.544"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.03%  "

$ws.Range("D30").Value = "'4.232"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "'4.187"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").Value = "'1.203"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").Value = "'0.05119"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.28%  "

$ws.Range("D34").Value = "'0.7962"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.28%  "

$ws.Range("D35").Value = "'1.917"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.93%  "

$ws.Range("D36").Value = "'1.148"
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("D38").Value = "'1.329.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.07%  "

$ws.Range("D39").Value = "'0.01874"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "

$ws.Range("D40").Value = "'2.717"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.9543"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.43%  "

$ws.Range("D42").Value = "'6.067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.64%  "

$ws.Range("D43").Value = "'107.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.10%  "

$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "'9.745"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("D47").Value = "'1.990.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "

$ws.Range("D48").Value = "'0.5181"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").Value = "'64.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.28%  "

$ws.Range("D50").Value = "'1.770"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.05%  "

$ws.Range("D51").Value = "'7.005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
